$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 8.403200000000005
$ws.Range("B12").Value = 6.3757
$ws.Range("D13").Value = -7.856200000000002
$ws.Range("B18").Value = 5.286700000000003
